$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts existing rows 9-81 down to 10-82)
$ws.Rows(9).Insert()

# Fill in the new row 9 with this week's data
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = 45163
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 100112012
$ws.Cells.Item(9, 7).Value = "Espinaca"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 35
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 13000
$ws.Cells.Item(9, 13).Value = 13000
$ws.Cells.Item(9, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(9, 15).Value = "Región Metropolitana"
$ws.Cells.Item(9, 16).Value = 1300
$ws.Cells.Item(9, 17).Value = 10
$ws.Cells.Item(9, 18).Value = "Hortaliza"
